{"js": "// Resume edit: add C#/F#/Python to \"Experienced\", add Julia (and re-shuffle) to\n// \"Familiar\", and split the \"scientific inquiry\" bullet into two bullets with\n// a new \"Advised leadership...\" bullet inserted before it.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\nlet experiencedPara = null;\nlet familiarPara = null;\nfor (const p of paragraphs.items) {\n  if (experiencedPara === null && p.text.indexOf(\"(Experienced)\") === 0) {\n    experiencedPara = p;\n  }\n  if (familiarPara === null && p.text.indexOf(\"(Familiar)\") === 0) {\n    familiarPara = p;\n  }\n}\n\n// ---------------------------------------------------------------------\n// 1) \"(Experienced)\" programming-languages bullet: append C#, F#, Python.\n// ---------------------------------------------------------------------\nif (experiencedPara) {\n  const expResults = experiencedPara.search(\n    \"PowerShell, ECMAScript (JavaScript), Node.js, HTML/XML/CSS\",\n    { matchCase: true }\n  );\n  expResults.load(\"items\");\n  await context.sync();\n  if (expResults.items.length > 0) {\n    expResults.items[0].insertText(\n      \"PowerShell, ECMAScript (JavaScript), Node.js, HTML/XML/CSS, C#, F#, Python\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n}\n\n// ---------------------------------------------------------------------\n// 2) \"(Familiar)\" programming-languages bullet: add Julia, drop C#/F#/Python\n//    (now listed under Experienced), and reorder the remaining items.\n// ---------------------------------------------------------------------\nif (familiarPara) {\n  const famResults = familiarPara.search(\n    \"C#, F#, ReasonML, Rust, WebAssembly, Kotlin, Python, R, Clojure, ClojureScript, Elm\",\n    { matchCase: true }\n  );\n  famResults.load(\"items\");\n  await context.sync();\n  if (famResults.items.length > 0) {\n    famResults.items[0].insertText(\n      \"Julia, ReasonML, Clojure, ClojureScript, Rust, WebAssembly, Kotlin, Elm, R\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n}\n\n// ---------------------------------------------------------------------\n// 3) Operations Research Analyst bullets: the old \"Utilize scientific\n//    inquiry...\" bullet becomes \"Advised leadership on emergent\n//    technologies (such as AI/ML) and methodologies.\" and a new bullet\n//    with the (now past-tense) \"Utilized scientific inquiry...\" text is\n//    inserted right after it.\n// ---------------------------------------------------------------------\nconst scientificInquiryText =\n  \"Utilize scientific inquiry in the independent development of \" +\n  \"mathematical models and computer programs to evaluate and predict \" +\n  \"the ability to support assigned projects, studies, or problems.\";\n\nconst siResults = body.search(scientificInquiryText, { matchCase: true });\nsiResults.load(\"items\");\nawait context.sync();\n\nif (siResults.items.length > 0) {\n  const siRange = siResults.items[0];\n  const siPara = siRange.paragraphs.getFirst();\n\n  // Insert the duplicated (past-tense) bullet right after the current one.\n  siPara.insertParagraph(\n    \"Utilized scientific inquiry in the independent development of \" +\n      \"mathematical models and computer programs to evaluate and predict \" +\n      \"the ability to support assigned projects, studies, or problems.\",\n    Word.InsertLocation.after\n  );\n  await context.sync();\n\n  // Replace the original bullet's text with the new \"Advised leadership\" copy.\n  siPara.insertText(\n    \"Advised leadership on emergent technologies (such as AI/ML) and methodologies.\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n", "ps1": "# Resume edit: add C#/F#/Python to \"Experienced\", add Julia (and re-shuffle)\n# to \"Familiar\", and split the \"scientific inquiry\" bullet into two bullets\n# with a new \"Advised leadership...\" bullet inserted before it.\n\n$doc = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1) \"(Experienced)\" programming-languages bullet: append C#, F#, Python.\n# ---------------------------------------------------------------------\n$count = $doc.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $para = $doc.Paragraphs($i)\n    if ($para.Range.Text.StartsWith(\"(Experienced)\")) {\n        $rng = $para.Range\n        $find = $rng.Find\n        $find.MatchCase = $true\n        $found = $find.Execute(\"PowerShell, ECMAScript (JavaScript), Node.js, HTML/XML/CSS\")\n        if ($found) {\n            $rng.Text = \"PowerShell, ECMAScript (JavaScript), Node.js, HTML/XML/CSS, C#, F#, Python\"\n        }\n        break\n    }\n}\n\n# ---------------------------------------------------------------------\n# 2) \"(Familiar)\" programming-languages bullet: add Julia, drop C#/F#/Python\n#    (now listed under Experienced), and reorder the remaining items.\n# ---------------------------------------------------------------------\n$count = $doc.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $para = $doc.Paragraphs($i)\n    if ($para.Range.Text.StartsWith(\"(Familiar)\")) {\n        $rng = $para.Range\n        $find = $rng.Find\n        $find.MatchCase = $true\n        $found = $find.Execute(\"C#, F#, ReasonML, Rust, WebAssembly, Kotlin, Python, R, Clojure, ClojureScript, Elm\")\n        if ($found) {\n            $rng.Text = \"Julia, ReasonML, Clojure, ClojureScript, Rust, WebAssembly, Kotlin, Elm, R\"\n        }\n        break\n    }\n}\n\n# ---------------------------------------------------------------------\n# 3) Operations Research Analyst bullets: the old \"Utilize scientific\n#    inquiry...\" bullet becomes \"Advised leadership on emergent\n#    technologies (such as AI/ML) and methodologies.\" and a new bullet\n#    with the (now past-tense) \"Utilized scientific inquiry...\" text is\n#    inserted right after it.\n# ---------------------------------------------------------------------\n$count = $doc.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $para = $doc.Paragraphs($i)\n    if ($para.Range.Text -like \"*scientific inquiry*\") {\n        $para.Range.InsertParagraphAfter()\n        $doc.Paragraphs($i + 1).Range.Text = \"Utilized scientific inquiry in the independent development of mathematical models and computer programs to evaluate and predict the ability to support assigned projects, studies, or problems.\"\n        $doc.Paragraphs($i).Range.Text = \"Advised leadership on emergent technologies (such as AI/ML) and methodologies.\"\n        break\n    }\n}\n"}
